# Andre_Master.xlsx: add "Lookup" column at the front of the header row
# and a "Lookup" table value below it (account-loading / lookup tables
# for account-to-ledgers per commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing headers one column to the right and insert the new
# "Lookup" header in column A. Snapshot the old values first (note: the
# getter must be invoked with explicit call parens, `.Value()`, otherwise
# it yields an unbound property descriptor instead of the cell's value).
$startBalance = $ws.Range("A1").Value()
$debitTotal = $ws.Range("B1").Value()
$creditTotal = $ws.Range("C1").Value()
$endBalance = $ws.Range("D1").Value()

$ws.Range("A1").Value = "Lookup"
$ws.Range("B1").Value = $startBalance
$ws.Range("C1").Value = $debitTotal
$ws.Range("D1").Value = $creditTotal
$ws.Range("E1").Value = $endBalance

# New lookup value under the new header, formatted with a 3-digit
# zero-padded custom number format.
$ws.Range("A2").Value = 1
$ws.Range("A2").NumberFormat = "000"

# Match the author's final selection position.
[void]$ws.Range("AD16").Select()
